# New profile script implementation
# Adds a new test case row (Profile46 / OPQA-2937) to the Profile test sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add row 47 ------------------------------------------------------
# Carry over the formatting used by the previous data row (row 46) so the
# new row matches the existing styles (borders, fonts, etc.)
$ws.Range("A46:E46").Copy() | Out-Null
$ws.Range("A47:E47").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A47").Value = "Profile46"
$ws.Range("B47").Value = "OPQA-2937 "
$ws.Range("C47").Value = "Verify that profile call to Action(CTA) in a white box is getting  displayed in country field when your first name, last name, title/role, institution, or location/country is blank"
$ws.Range("D47").Value = "Y"

# --- Column C got wider to accommodate the new (longer) description --
$ws.Columns.Item(3).ColumnWidth = 173

# --- Sheet view: scroll down and move the active selection -----------
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("C31").Select() | Out-Null
